$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the cryptos list refresh (diff-derived).
$ws.Range("D2").Value = "56.114.88"
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("D3").Value = "2.975.07"
$ws.Range("E3").Value = "  -4.97%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'495.31"
$ws.Range("E5").Value = "  -5.36%  "
$ws.Range("D6").Value = "'134.43"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "2.969.37"
$ws.Range("E8").Value = "  -5.07%  "
$ws.Range("D9").Value = "'0.424"
$ws.Range("E9").Value = "  -3.37%  "
$ws.Range("D10").Value = "'7.24"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "  -6.84%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "3.491.73"
$ws.Range("E14").Value = "  -5.01%  "
$ws.Range("D15").Value = "'24.84"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "56.204.58"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.988.82"
$ws.Range("E17").Value = "  -4.84%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000145"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("D19").Value = "'5.79"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'12.31"
$ws.Range("E20").Value = "  -5.61%  "
$ws.Range("D21").Value = "'7.70"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "'323.97"
$ws.Range("E22").Value = "  -5.38%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'0.463"
$ws.Range("E24").Value = "  -8.29%  "
$ws.Range("D25").Value = "'61.00"
$ws.Range("E25").Value = "  -10.11%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("D28").Value = "0.0₃0894"
$ws.Range("E28").Value = "  -5.57%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'6.44"
$ws.Range("E30").Value = "  -4.85%  "
$ws.Range("D31").Value = "'6.74"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("D33").Value = "'1.72"
$ws.Range("E33").Value = "  -7.20%  "
$ws.Range("D34").Value = "'19.95"
$ws.Range("E34").Value = "  -6.79%  "
$ws.Range("D35").Value = "'154.64"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "'4.46"
$ws.Range("E36").Value = "  -6.73%  "
$ws.Range("D37").Value = "'1.27"
$ws.Range("E37").Value = "  -6.91%  "
$ws.Range("D38").Value = "'5.57"
$ws.Range("E38").Value = "  -10.07%  "
$ws.Range("D39").Value = "'0.0674"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("D40").Value = "'23.18"
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("D41").Value = "3.013.31"
$ws.Range("E41").Value = "  -4.81%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'36.17"
$ws.Range("E43").Value = "  -10.10%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -7.18%  "
$ws.Range("D45").Value = "'0.634"
$ws.Range("E45").Value = "  -7.92%  "
$ws.Range("D46").Value = "'1.40"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").Value = "2.190.55"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").Value = "'3.54"
$ws.Range("E48").Value = "  -8.90%  "
$ws.Range("D49").Value = "'1.93"
$ws.Range("E49").Value = "  +5.92%  "
$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("E51").Value = "  -7.53%  "

Write-Host "Applied cryptos list update"
